$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary row: average of column J (the k value) in row 12
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary label + value rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the first summary value: bold, size 12, vertically centered
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

# Copy that formatting onto the other three summary values (format painter)
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select the summary block, as the author last left it selected
$ws.Range("A14:B17").Select()

# Page setup: portrait, paper size 9 (A4)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
